$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look numeric need NumberFormat "@" (Text)
# applied first so Excel does not silently convert them to numbers,
# which would also rstrip significant trailing zeros (e.g. "70.20" -> 70.2).

$ws.Range('D2').Value = '37.054.44'
$ws.Range('E2').Value = '  +1.33%  '

$ws.Range('D3').Value = '1.985.13'
$ws.Range('E3').Value = '  +1.18%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.59'
$ws.Range('E5').Value = '  +0.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('E6').Value = '  +1.96%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.97'
$ws.Range('E7').Value = '  +3.17%  '

$ws.Range('E9').Value = '  +2.22%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0801'
$ws.Range('E10').Value = '  -1.38%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +0.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.05'
$ws.Range('E12').Value = '  +9.79%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.21'
$ws.Range('E13').Value = '  -0.17%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.846'
$ws.Range('E14').Value = '  +2.37%  '

$ws.Range('D15').Value = '2.276.26'
$ws.Range('E15').Value = '  +1.14%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.48'
$ws.Range('E16').Value = '  +4.02%  '

$ws.Range('D17').Value = '1.985.80'
$ws.Range('E17').Value = '  +1.55%  '

$ws.Range('D18').Value = '36.885.72'
$ws.Range('E18').Value = '  +1.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.20'
$ws.Range('E19').Value = '  +0.41%  '

$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  +0.45%  '

$ws.Range('E21').Value = '  +2.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.05'
$ws.Range('E22').Value = '  +0.49%  '

$ws.Range('E23').Value = '  +0.03%  '

$ws.Range('E24').Value = '  +2.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +0.12%  '

$ws.Range('E26').Value = '  +6.75%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  +0.94%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.51'
$ws.Range('E28').Value = '  +2.12%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.55'
$ws.Range('E29').Value = '  +0.53%  '

$ws.Range('E30').Value = '  +17.95%  '

$ws.Range('E31').Value = '  +2.13%  '

$ws.Range('E32').Value = '  +3.30%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0622'
$ws.Range('E33').Value = '  +0.42%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.53'
$ws.Range('E34').Value = '  +6.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.31'
$ws.Range('E35').Value = '  +2.56%  '

$ws.Range('E36').Value = '  +0.04%  '

$ws.Range('E37').Value = '  +0.02%  '

$ws.Range('E38').Value = '  +0.29%  '

$ws.Range('E39').Value = '  -7.04%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0988'
$ws.Range('E40').Value = '  +0.37%  '

$ws.Range('E41').Value = '  +0.98%  '

$ws.Range('E42').Value = '  +0.83%  '

$ws.Range('E43').Value = '  +0.99%  '

$ws.Range('E44').Value = '  +3.44%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.20'
$ws.Range('E45').Value = '  +2.86%  '

$ws.Range('D46').Value = '1.370.08'
$ws.Range('E46').Value = '  +0.39%  '

$ws.Range('E47').Value = '  +0.31%  '

$ws.Range('E48').Value = '  +1.78%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.49'
$ws.Range('E49').Value = '  +6.29%  '

$ws.Range('E50').Value = '  -0.44%  '

$ws.Range('E51').Value = '  +10.19%  '
